$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210..321 down to 211..322
$ws.Rows(210).Insert()

# Populate the newly inserted row 210 with its values
$ws.Range("A210").Value = 10
$ws.Range("B210").Value = "Vega Modelo de Temuco"
$ws.Range("C210").Value = "La Araucanía"
$ws.Range("D210").Value = 44488
$ws.Range("E210").Value = 9
$ws.Range("F210").Value = 100112043
$ws.Range("G210").Value = "Pepino ensalada"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 125
$ws.Range("K210").Value = 12000
$ws.Range("L210").Value = 12000
$ws.Range("M210").Value = 12000
$ws.Range("N210").Value = "$/caja 60 unidades"
$ws.Range("O210").Value = "Región de Arica y Parinacota"
$ws.Range("P210").Value = 200
$ws.Range("Q210").Value = 60
$ws.Range("R210").Value = "Hortaliza"
